$wb = $excel.ActiveWorkbook

# Update the per-capita cost driver (D14) on both scenario sheets.
# Value increases tenfold: 38377 -> 383775. All dependent formulas
# (N3, N4, N6, N13:N41, etc.) recalculate automatically.
$ws1 = $wb.Worksheets.Item("Test1")
$ws1.Range("D14").Value = 383775

$ws2 = $wb.Worksheets.Item("Test 2")
$ws2.Range("D14").Value = 383775
